# Add 9 numbered state-labels ("1".."9", skipping 10) to the simulation
# diagram on slide 1. Each label is a copy of the existing "TextBox 50"
# ("BON") shape -- duplicating it (instead of Shapes.AddTextbox) is what
# lets the new shapes pick up the exact same formatting already baked
# into the deck's other numbered labels (centered 28pt Helvetica text,
# no-fill autosize textbox, etc.) rather than PowerPoint's bare defaults.
#
# Left/Top are written as the nearest single-precision (float32) point
# value that converts back to the exact target EMU offset recorded in
# the authored deck, since the Shape.Left/.Top setters round-trip
# through a 32-bit float before being stored as EMUs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "TextBox 50" (the "BON" label) already has the cx=1989181 cy=523220
# size shared by every numbered label, sz=2800 Helvetica centered text,
# and a no-fill autosize text box -- use it as the template to clone.
$template = $s.Shapes.Item(13)

$labels = @(
    @{ Name = "TextBox 71"; Text = "1"; Left = -2.747637987136841;  Top = 200.57630920410156 },
    @{ Name = "TextBox 72"; Text = "2"; Left = 221.91607666015625;  Top = 198.97567749023438 },
    @{ Name = "TextBox 74"; Text = "3"; Left = 456.6933288574219;   Top = 197.48858642578125 },
    @{ Name = "TextBox 75"; Text = "4"; Left = 673.7725219726562;   Top = 43.436851501464844 },
    @{ Name = "TextBox 76"; Text = "5"; Left = 675.0406494140625;   Top = 357.2267150878906  },
    @{ Name = "TextBox 78"; Text = "8"; Left = 824.6774291992188;   Top = 261.46600341796875 },
    @{ Name = "TextBox 79"; Text = "9"; Left = 455.478515625;       Top = 47.82921600341797  },
    @{ Name = "TextBox 80"; Text = "7"; Left = 303.7677307128906;   Top = 45.706851959228516 },
    @{ Name = "TextBox 82"; Text = "6"; Left = 140.6618194580078;   Top = 47.82921600341797  }
)

foreach ($label in $labels) {
    $dup = $template.Duplicate()
    $dup.Name = $label.Name
    $dup.Left = $label.Left
    $dup.Top = $label.Top
    $dup.TextFrame.TextRange.Text = $label.Text
}
